# Rename sheets and a couple columns for consistency

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "MarketParams"

# Update header text in column A (was "Average price", now "year")
$ws.Range("A1").Value = "year"

# Update the active selection to match the target state
$ws.Range("C5").Select()
